$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.9020979004458534, 1.2629696074199597, -0.5052416671458982, 2.7667445119319796)"
$ws.Range("C2").Value = "NIG(1.0918022350261651, 0.9040073947564247, 3.3551220385533287, 5.684550655994077)"
$ws.Range("D2").Value = "NCT(3.431552188186929, 1.5072548954658505, -0.6637997885330891, 2.296524586821742)"
$ws.Range("E2").Value = "NIG(1.0529047686233677, 0.7107511552118035, 5.798773274556774, 5.320410267797112)"
